# Apply updated FFXIV market-board pricing figures to each job sheet's
# profit table (currentAveragePrice* / LevePrice* / LeveProfit* columns).
# Values below are the refreshed averages pulled by the scheduled runner;
# only numeric H:N cells change (plus one newly-populated HQ profit cell).
$wb = $excel.ActiveWorkbook

# --- ALC sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
# Row 2: H2, I2, K2, M2
$ws.Range("H2").Value = 222.5
$ws.Range("I2").Value = 222.5
$ws.Range("K2").Value = 222.5
$ws.Range("M2").Value = -109.5
# Row 9: H9, I9, J9, K9, L9, M9, N9
$ws.Range("H9").Value = 101.5
$ws.Range("I9").Value = 111.166664
$ws.Range("J9").Value = 72.5
$ws.Range("K9").Value = 111.166664
$ws.Range("L9").Value = 72.5
$ws.Range("M9").Value = 57.833336
$ws.Range("N9").Value = -410.5
# Row 15: H15, I15, K15, M15
$ws.Range("H15").Value = 1990.2162
$ws.Range("I15").Value = 1990.2162
$ws.Range("K15").Value = 5970.6486
$ws.Range("M15").Value = -5801.6486
# Row 40: H40, I40, J40, K40, L40, M40, N40
$ws.Range("H40").Value = 5930.8887
$ws.Range("I40").Value = 2750
$ws.Range("J40").Value = 8475.6
$ws.Range("K40").Value = 2750
$ws.Range("L40").Value = 8475.6
$ws.Range("M40").Value = -2575
$ws.Range("N40").Value = -8825.6
# Row 41: H41, I41, K41, M41
$ws.Range("H41").Value = 6305.8423
$ws.Range("I41").Value = 751.4545000000001
$ws.Range("K41").Value = 751.4545000000001
$ws.Range("M41").Value = -311.4545000000001
# Row 43: H43, I43, K43, M43
$ws.Range("H43").Value = 2079.1765
$ws.Range("I43").Value = 1639.5834
$ws.Range("K43").Value = 1639.5834
$ws.Range("M43").Value = -1570.5834
# Row 59: H59, I59, J59, K59, L59, M59, N59
$ws.Range("H59").Value = 300
$ws.Range("I59").Value = 300
$ws.Range("J59").Value = 300
$ws.Range("K59").Value = 900
$ws.Range("L59").Value = 900
$ws.Range("M59").Value = -343
$ws.Range("N59").Value = -2014
# Row 64: H64, I64, J64, K64, L64, M64, N64
$ws.Range("H64").Value = 7061.231
$ws.Range("I64").Value = 5690.4614
$ws.Range("J64").Value = 8432
$ws.Range("K64").Value = 5690.4614
$ws.Range("L64").Value = 8432
$ws.Range("M64").Value = -5442.4614
$ws.Range("N64").Value = -8928
# Row 67: H67, I67, J67, K67, L67, M67, N67
$ws.Range("H67").Value = 7061.231
$ws.Range("I67").Value = 5690.4614
$ws.Range("J67").Value = 8432
$ws.Range("K67").Value = 5690.4614
$ws.Range("L67").Value = 8432
$ws.Range("M67").Value = -4832.4614
$ws.Range("N67").Value = -10148
# Row 74: H74, I74, K74, M74
$ws.Range("H74").Value = 6711.7334
$ws.Range("I74").Value = 4606.909
$ws.Range("K74").Value = 4606.909
$ws.Range("M74").Value = -3670.909
# Row 77: H77, I77, K77, M77
$ws.Range("H77").Value = 6711.7334
$ws.Range("I77").Value = 4606.909
$ws.Range("K77").Value = 23034.545
$ws.Range("M77").Value = -18354.545
# Row 86: H86, I86, K86, M86
$ws.Range("H86").Value = 3679.6428
$ws.Range("I86").Value = 3960.111
$ws.Range("K86").Value = 3960.111
$ws.Range("M86").Value = -2837.111
# Row 89: H89, I89, K89, M89
$ws.Range("H89").Value = 3679.6428
$ws.Range("I89").Value = 3960.111
$ws.Range("K89").Value = 19800.555
$ws.Range("M89").Value = -14184.555
# Row 106: H106, I106, K106, M106
$ws.Range("H106").Value = 2862.8
$ws.Range("I106").Value = 2862.8
$ws.Range("K106").Value = 2862.8
$ws.Range("M106").Value = -2231.8
# Row 112: H112, J112, L112, N112
$ws.Range("H112").Value = 1882.7317
$ws.Range("J112").Value = 1962.2778
$ws.Range("L112").Value = 5886.8334
$ws.Range("N112").Value = -8102.8334
# Row 135: H135, J135, L135, N135
$ws.Range("H135").Value = 17247444
$ws.Range("J135").Value = 16805
$ws.Range("L135").Value = 151245
$ws.Range("N135").Value = -156315
# Row 138: H138, I138, J138, K138, L138, M138, N138
$ws.Range("H138").Value = 1793.75
$ws.Range("I138").Value = 1181.3462
$ws.Range("J138").Value = 2931.0715
$ws.Range("K138").Value = 3544.0386
$ws.Range("L138").Value = 8793.2145
$ws.Range("M138").Value = 1595.9614
$ws.Range("N138").Value = -19073.2145

# --- ARM sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
# Row 63: H63, J63, L63, N63
$ws.Range("H63").Value = 4215.143
$ws.Range("J63").Value = 4453
$ws.Range("L63").Value = 4453
$ws.Range("N63").Value = -5825
# Row 66: H66, J66, L66, N66
$ws.Range("H66").Value = 4215.143
$ws.Range("J66").Value = 4453
$ws.Range("L66").Value = 22265
$ws.Range("N66").Value = -29129

# --- BSM sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
# Row 20: H20, I20, K20, M20
$ws.Range("H20").Value = 1900.1666
$ws.Range("I20").Value = 1973.3334
$ws.Range("K20").Value = 1973.3334
$ws.Range("M20").Value = -1726.3334
# Row 86: H86, I86, K86, M86
$ws.Range("H86").Value = 2760.8823
$ws.Range("I86").Value = 1561.1765
$ws.Range("K86").Value = 1561.1765
$ws.Range("M86").Value = -438.1765
# Row 89: H89, I89, K89, M89
$ws.Range("H89").Value = 2760.8823
$ws.Range("I89").Value = 1561.1765
$ws.Range("K89").Value = 7805.8825
$ws.Range("M89").Value = -2189.8825

# --- CRP sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31, I31, K31, M31
$ws.Range("H31").Value = 3389.32
$ws.Range("I31").Value = 1338.3077
$ws.Range("K31").Value = 1338.3077
$ws.Range("M31").Value = -1043.3077
# Row 34: H34, I34, K34, M34
$ws.Range("H34").Value = 3389.32
$ws.Range("I34").Value = 1338.3077
$ws.Range("K34").Value = 1338.3077
$ws.Range("M34").Value = -1136.3077
# Row 58: H58, I58, K58, M58
$ws.Range("H58").Value = 1335.8148
$ws.Range("I58").Value = 1034.68
$ws.Range("K58").Value = 1034.68
$ws.Range("M58").Value = -831.6800000000001
# Row 59: H59, J59, L59, N59
$ws.Range("H59").Value = 235000
$ws.Range("J59").Value = 70000
$ws.Range("L59").Value = 70000
$ws.Range("N59").Value = -72290
# Row 99: H99, J99, L99, N99
$ws.Range("H99").Value = 6606376
$ws.Range("J99").Value = 10007078
$ws.Range("L99").Value = 10007078
$ws.Range("N99").Value = -10010074
# Row 122: H122, I122, K122, M122
$ws.Range("H122").Value = 331125.97
$ws.Range("I122").Value = 394090.62
$ws.Range("K122").Value = 1182271.86
$ws.Range("M122").Value = -1179821.86
# Row 126: H126, J126, L126, N126
$ws.Range("H126").Value = 6606376
$ws.Range("J126").Value = 10007078
$ws.Range("L126").Value = 30021234
$ws.Range("N126").Value = -30026174
# Row 136: H136, I136, K136, M136
$ws.Range("H136").Value = 1335.8148
$ws.Range("I136").Value = 1034.68
$ws.Range("K136").Value = 3104.04
$ws.Range("M136").Value = -554.04

# --- CUL sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
# Row 92: H92, J92, L92, N92
$ws.Range("H92").Value = 600
$ws.Range("J92").Value = 705.5
$ws.Range("L92").Value = 2116.5
$ws.Range("N92").Value = -4612.5
# Row 113: H113, I113, J113, K113, L113, M113, N113
$ws.Range("H113").Value = 4802.115
$ws.Range("I113").Value = 8486.691999999999
$ws.Range("J113").Value = 1117.5385
$ws.Range("K113").Value = 25460.076
$ws.Range("L113").Value = 3352.6155
$ws.Range("M113").Value = -23290.076
$ws.Range("N113").Value = -7692.6155

# --- GSM sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
# Row 70: H70, I70, J70, K70, L70, M70, N70
$ws.Range("H70").Value = 107086.45
$ws.Range("I70").Value = 226303.2
$ws.Range("J70").Value = 7739.1665
$ws.Range("K70").Value = 226303.2
$ws.Range("L70").Value = 7739.1665
$ws.Range("M70").Value = -226033.2
$ws.Range("N70").Value = -8279.166499999999
# Row 73: H73, I73, J73, K73, L73, M73, N73
$ws.Range("H73").Value = 107086.45
$ws.Range("I73").Value = 226303.2
$ws.Range("J73").Value = 7739.1665
$ws.Range("K73").Value = 226303.2
$ws.Range("L73").Value = 7739.1665
$ws.Range("M73").Value = -225367.2
$ws.Range("N73").Value = -9611.166499999999
# Row 97: H97, I97, K97, M97
$ws.Range("H97").Value = 408.07407
$ws.Range("I97").Value = 328.88235
$ws.Range("K97").Value = 328.88235
$ws.Range("M97").Value = 167.11765
# Row 122: H122, I122, K122, M122
$ws.Range("H122").Value = 5857.25
$ws.Range("I122").Value = 5628.7
$ws.Range("K122").Value = 16886.1
$ws.Range("M122").Value = -14436.1
# Row 132: H132, I132, K132, M132
$ws.Range("H132").Value = 3206.3845
$ws.Range("I132").Value = 2878.2727
$ws.Range("K132").Value = 8634.8181
$ws.Range("M132").Value = -6104.8181

# --- LTW sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
# Row 46: H46, I46, K46, M46
$ws.Range("H46").Value = 8544.5
$ws.Range("I46").Value = 1938.5
$ws.Range("K46").Value = 1938.5
$ws.Range("M46").Value = -1750.5
# Row 61: H61, J61, L61, N61
$ws.Range("H61").Value = 1986.65
$ws.Range("J61").Value = 7000
$ws.Range("L61").Value = 7000
$ws.Range("N61").Value = -7404
# Row 68: H68, I68, J68, K68, L68, M68, N68
$ws.Range("H68").Value = 4633.6665
$ws.Range("I68").Value = 3284.3572
$ws.Range("J68").Value = 5492.3184
$ws.Range("K68").Value = 3284.3572
$ws.Range("L68").Value = 5492.3184
$ws.Range("M68").Value = -2535.3572
$ws.Range("N68").Value = -6990.3184
# Row 71: H71, I71, J71, K71, L71, M71, N71
$ws.Range("H71").Value = 4633.6665
$ws.Range("I71").Value = 3284.3572
$ws.Range("J71").Value = 5492.3184
$ws.Range("K71").Value = 16421.786
$ws.Range("L71").Value = 27461.592
$ws.Range("M71").Value = -12677.786
$ws.Range("N71").Value = -34949.592
# Row 82: H82, I82, K82, M82
$ws.Range("H82").Value = 7279
$ws.Range("I82").Value = 3749.75
$ws.Range("K82").Value = 3749.75
$ws.Range("M82").Value = -3388.75
# Row 85: H85, I85, K85, M85
$ws.Range("H85").Value = 7279
$ws.Range("I85").Value = 3749.75
$ws.Range("K85").Value = 3749.75
$ws.Range("M85").Value = -2501.75
# Row 113: H113, J113, L113, N113
$ws.Range("H113").Value = 1986.65
$ws.Range("J113").Value = 7000
$ws.Range("L113").Value = 7000
$ws.Range("N113").Value = -11340
# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 4427.074
$ws.Range("I132").Value = 3584.05
$ws.Range("J132").Value = 6835.7144
$ws.Range("K132").Value = 10752.15
$ws.Range("L132").Value = 20507.1432
$ws.Range("M132").Value = -8222.150000000001
$ws.Range("N132").Value = -25567.1432

# --- WVR sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
# Row 107: H107, I107, J107, K107, L107, M107, N107
$ws.Range("H107").Value = 242.5
$ws.Range("I107").Value = 185
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 555
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 1365
$ws.Range("N107").Value = -4740
# Row 113: H113, I113, J113, K113, L113, M113, N113
$ws.Range("H113").Value = 445.44446
$ws.Range("I113").Value = 341.125
$ws.Range("J113").Value = 1280
$ws.Range("K113").Value = 1023.375
$ws.Range("L113").Value = 3840
$ws.Range("M113").Value = 1146.625
$ws.Range("N113").Value = -8180
# Row 122: H122, I122, K122, M122
$ws.Range("H122").Value = 6515.7085
$ws.Range("I122").Value = 3520.6667
$ws.Range("K122").Value = 10562.0001
$ws.Range("M122").Value = -8112.000100000001
# Row 133: H133, J133, L133, N133
$ws.Range("H133").Value = 30000
$ws.Range("J133").Value = 30000
$ws.Range("L133").Value = 30000
$ws.Range("N133").Value = -40120
